# RPA datasets push 2024-06-12
# Insert two new IPO entries (키움스팩9호, 넥스트바이오메디컬) right after the
# first data row, shifting every subsequent row down by two. Because the
# table always keeps exactly 20 data rows (A1:F21), the two oldest entries
# that used to sit at the bottom (신한글로벌액티브, 한중엔시에스) drop off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the old row 3 (after 아이빔테크놀로지).
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# New row 3: 키움스팩9호
$ws.Range("A3").Value = "키움스팩9호"
$ws.Range("B3").Value = "2024.07.02~07.03"
$ws.Range("C3").Value = "2,000~2,000"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 6000
$ws.Range("F3").Value = "키움증권"

# New row 4: 넥스트바이오메디컬
$ws.Range("A4").Value = "넥스트바이오메디컬"
$ws.Range("B4").Value = "2024.07.01~07.05"
$ws.Range("C4").Value = "24,000~29,000"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = 24000
$ws.Range("F4").Value = "한국투자증권"

# 미래에셋비전스팩6호's 확정공모가 (final offering price) was finalized;
# it moved from old row 13 to row 15 after the insert above. Copy the text
# value "2000" from a sibling cell (D16) so it is stored as text, not a
# number, matching the rest of column D.
$ws.Range("D16").Copy()
$ws.Range("D15").PasteSpecial()
$excel.CutCopyMode = $false

# Drop the two oldest rows that fell off the bottom of the fixed-size table
# (now sitting at rows 22 and 23 after the insert above).
$ws.Rows.Item(22).EntireRow.Delete()
$ws.Rows.Item(22).EntireRow.Delete()
